$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f = $sec.Footers(2)  # footer1.xml id=2
$ishp = $f.Range.InlineShapes(1)
$r = $ishp.Range
$xml = $r.WordOpenXML
if ($xml -match '(?s)<w:drawing>.*?</w:drawing>') {
    $drawing = $matches[0]
}
$newDrawing = $drawing -replace 'name="image2\.png"', 'name="image1.png"'

# Delete the shape first (keeps rId1 relationship entry in rels, unused).
$r.Select() | Out-Null
$ishp.Delete()

$minimalPkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:r>' + $newDrawing + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($minimalPkg)
